# Logged Week 16 and performed season sim from Week 17
$wb = $excel.ActiveWorkbook

# --- OFF sheet: update Row "R" (row 3) with new cumulative totals ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 215
$wsOff.Range("C3").Value = 149
$wsOff.Range("D3").Value = 51
$wsOff.Range("E3").Value = 22
$wsOff.Range("F3").Value = 8

# --- DEF sheet: update Row "R" (row 3) with new cumulative totals ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 199
$wsDef.Range("C3").Value = 132
$wsDef.Range("D3").Value = 55
$wsDef.Range("E3").Value = 30
